$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (values + copy style/format from the neighboring header cell)
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Updated B/C/D values, and new G/H values per row
$data = @{
  2  = @{ B = 0.5012650052568042;  C = 0.990018749248315;  D = 0.5783909156722671 }
  3  = @{ B = 0.2408662486245789;  C = 0.9952922631133009; D = 0.3906853379118294 }
  4  = @{ B = 0.2796426285527719;  C = 0.9946202976515992; D = 0.4280836761227448 }
  5  = @{ B = 0.4107685786646339;  C = 0.9919002060210602; D = 0.4918168826609312 }
  6  = @{ B = 0.4838554589956822;  C = 0.985787897678873;  D = 0.5111291360470315 }
  7  = @{ B = 0.09370646054035536; C = 0.9987049513429616; D = 0.2492622490976881 }
  8  = @{ B = 0.03283384250914684; C = 0.9996615885863042; D = 0.1307150721627474 }
  9  = @{ B = 0.09490897099941178; C = 0.999434654339443;  D = 0.2166318052638244 }
  10 = @{ B = 0.06675309232430304; C = 0.9987805790882839; D = 0.2050030884974856 }
  11 = @{ B = 0.1246800155313699;  C = 0.9907881862466283; D = 0.2717651402261659 }
  12 = @{ B = 0.05249095752039892; C = 0.9984563304914589; D = 0.1675798128709076 }
  13 = @{ B = 0.05951985533614699; C = 0.9994354274984324; D = 0.1779264057377217 }
  14 = @{ B = 0.05599300176264916; C = 0.9992418893195172; D = 0.1915715355941557 }
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  $ws.Range("B$row").Value = $vals.B
  $ws.Range("C$row").Value = $vals.C
  $ws.Range("D$row").Value = $vals.D
  $ws.Range("G$row").Value = 1.669922641383406
  $ws.Range("H$row").Value = 0.97
}
